$d = $word.ActiveDocument

function Split-LeftRun([string]$oldWord) {
    # Locate the run containing e.g. "left-top" / "left-bottom" and split it
    # into "left-" (keeps the original run formatting) + "middle" (plain
    # run, formatting copied from a donor range that already carries no
    # explicit color - just like Word does when it creates a fresh run).
    $rng = $d.Content
    $rng.Find.Execute($oldWord) | Out-Null

    $prefixEnd = $rng.Start + 5  # length of "left-"
    $suffixRange = $d.Range($prefixEnd, $rng.End)
    $suffixRange.Text = ""   # delete "top"/"bottom", leaving "left-"

    $insertionPoint = $d.Range($prefixEnd, $prefixEnd)
    $insertionPoint.InsertAfter("middle")

    $newRunRange = $d.Range($prefixEnd, $prefixEnd + 6)

    # Donor range: an existing plain-formatted run (rPr with only <w:rtl/>,
    # no <w:color>) elsewhere in the document, used so the freshly created
    # run ends up with that same bare formatting instead of inheriting the
    # colored formatting of its neighbor.
    $donorRng = $d.Content
    $donorRng.Find.Execute("<env>a") | Out-Null
    $donor = $d.Range($donorRng.End - 1, $donorRng.End)

    $newRunRange.FormattedText = $donor.FormattedText
    $midRange = $d.Range($prefixEnd, $prefixEnd + 1)
    $midRange.Text = "middle"
}

Split-LeftRun "left-top"
Split-LeftRun "left-bottom"

$d.PageSetup.FooterDistance = 36
